$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at 16 (one new "Clean Code" log entry), formatted like
#    the existing data rows (copy formatting from row 15).
# ---------------------------------------------------------------------------
$ws.Rows(16).Insert()
$ws.Range("A15:J15").Copy()
$ws.Range("A16:J16").PasteSpecial(-4122)
$ws.Rows(16).RowHeight = 15

# ---------------------------------------------------------------------------
# 2) Insert two new rows at 19-20 (two more new "Clean Code" log entries).
# ---------------------------------------------------------------------------
$ws.Rows("19:20").Insert()
$ws.Range("A17:J17").Copy()
$ws.Range("A19:J20").PasteSpecial(-4122)
$ws.Rows(19).RowHeight = 15
$ws.Rows(20).RowHeight = 15

# ---------------------------------------------------------------------------
# 3) Column A: drop the helper "=prev+1" formulas, replace with literal
#    sequential numbers 1..20.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11
$ws.Range("A18").Value = 12
$ws.Range("A19").Value = 13
$ws.Range("A20").Value = 14
$ws.Range("A21").Value = 15
$ws.Range("A22").Value = 16
$ws.Range("A23").Value = 17
$ws.Range("A24").Value = 18
$ws.Range("A25").Value = 19
$ws.Range("A26").Value = 20

# ---------------------------------------------------------------------------
# 4) Fill in the data for the three brand-new "Clean Code" rows (16, 19, 20).
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 0.95486111111111116
$ws.Range("D16").Value = 0.96527777777777779
$ws.Range("E16").Value = 0
$ws.Range("F16").Formula = "=(D16-C16)*1440"
$ws.Range("G16").Value = "Clean Code"

$ws.Range("C19").Value = 0.73055555555555562
$ws.Range("D19").Value = 0.75347222222222221
$ws.Range("E19").Value = 0
$ws.Range("F19").Formula = "=(D19-C19)*1440"
$ws.Range("G19").Value = "Clean Code"

$ws.Range("C20").Value = 0.95208333333333339
$ws.Range("D20").Value = 0.97499999999999998
$ws.Range("E20").Value = 0
$ws.Range("F20").Formula = "=(D20-C20)*1440"
$ws.Range("G20").Value = "Clean Code"

# ---------------------------------------------------------------------------
# 5) Fill in the previously-blank row (now row 24) with the new
#    "VL27 / LambdaExpression" log entry.
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 0.75
$ws.Range("D24").Value = 0.76666666666666661
$ws.Range("E24").Value = 0
$ws.Range("G24").Value = "VL27"
$ws.Range("H24").Value = "LambdaExpression"
$ws.Range("I24").Value = "x"

# ---------------------------------------------------------------------------
# 6) View / selection bookkeeping to match the final workbook state.
# ---------------------------------------------------------------------------
$ws.Range("L21").Select()
$excel.ActiveWindow.ScrollRow = 1

Write-Host "edit complete"
